# Update countries & provincias Spain
# Applies the refreshed COVID country-stats snapshot (6 Abril 2020, 17:52)
# to the "Pais" worksheet: a handful of countries swap positions (because
# the sheet is kept sorted by "Casos totales") and their numeric columns
# (B:H) are refreshed to the newer figures, plus the "last updated" banner
# in A1 gets its new timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Banner timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 17:52"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 339131
$ws.Range("C4").Value = 2458
$ws.Range("D4").Value = 18029
$ws.Range("E4").Value = 311413
$ws.Range("G4").Value = 73
$ws.Range("H4").Value = 9689

# --- Row 5: España ---
$ws.Range("E5").Value = 81426
$ws.Range("G5").Value = 528
$ws.Range("H5").Value = 13169

# --- Row 7: Alemania ---
$ws.Range("B7").Value = 100920
$ws.Range("C7").Value = 797
$ws.Range("E7").Value = 70611
$ws.Range("G7").Value = 25
$ws.Range("H7").Value = 1609

# --- Row 11: Reino Unido ---
$ws.Range("B11").Value = 51608
$ws.Range("C11").Value = 3802
$ws.Range("E11").Value = 46100
$ws.Range("G11").Value = 439
$ws.Range("H11").Value = 5373

# --- Row 16: Canada ---
$ws.Range("B16").Value = 15853
$ws.Range("C16").Value = 341
$ws.Range("D16").Value = 3128
$ws.Range("E16").Value = 12432
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 293

# --- Rows 28-29: Dinamarca / Chequia swap position ---
$ws.Range("A28").Value = "Chequia"
$ws.Range("B28").Value = 4735
$ws.Range("C28").Value = 148
$ws.Range("D28").Value = 121
$ws.Range("E28").Value = 4536
$ws.Range("F28").Value = 84
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = 78

$ws.Range("A29").Value = "Dinamarca"
$ws.Range("B29").Value = 4681
$ws.Range("C29").Value = 312
$ws.Range("D29").Value = 1378
$ws.Range("E29").Value = 3116
$ws.Range("F29").Value = 144
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = 187

# --- Row 31: Polonia ---
$ws.Range("B31").Value = 4413
$ws.Range("C31").Value = 311
$ws.Range("E31").Value = 4144
$ws.Range("G31").Value = 13
$ws.Range("H31").Value = 107

# --- Row 32: Rumania ---
$ws.Range("E32").Value = 3483
$ws.Range("G32").Value = 17
$ws.Range("H32").Value = 168

# --- Row 38: Luxemburgo ---
$ws.Range("B38").Value = 2843
$ws.Range("C38").Value = 39
$ws.Range("E38").Value = 2302
$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 41

# --- Rows 55-57: Argelia / Ucrania / Singapur reshuffle ---
$ws.Range("A55").Value = "Singapur"
$ws.Range("B55").Value = 1375
$ws.Range("C55").Value = 66
$ws.Range("D55").Value = 344
$ws.Range("E55").Value = 1025
$ws.Range("F55").Value = 25
$ws.Range("H55").Value = 6

$ws.Range("A56").Value = "Argelia"
$ws.Range("B56").Value = 1320
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 90
$ws.Range("E56").Value = 1078
$ws.Range("F56").Value = 46
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 152

$ws.Range("A57").Value = "Ucrania"
$ws.Range("B57").Value = 1319
$ws.Range("C57").Value = 11
$ws.Range("D57").Value = 28
$ws.Range("E57").Value = 1253
$ws.Range("F57").Value = 16
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 38

# --- Row 88: Uzbekistan ---
$ws.Range("B88").Value = 397
$ws.Range("C88").Value = 55
$ws.Range("E88").Value = 365

# --- Rows 92-96: Burkina Faso / Jordania / Reunion / Oman / Cuba reshuffle ---
$ws.Range("A92").Value = "Cuba"
$ws.Range("B92").Value = 350
$ws.Range("C92").Value = 30
$ws.Range("D92").Value = 18
$ws.Range("E92").Value = 323
$ws.Range("F92").Value = 11
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 9

$ws.Range("A93").Value = "Reunion"
$ws.Range("B93").Value = 349
$ws.Range("C93").Value = 5
$ws.Range("D93").Value = 40
$ws.Range("E93").Value = 309
$ws.Range("F93").Value = 4
$ws.Range("H93").Value = 0

$ws.Range("A94").Value = "Burkina Faso"
$ws.Range("B94").Value = 345
$ws.Range("D94").Value = 90
$ws.Range("E94").Value = 238
$ws.Range("F94").Value = 0
$ws.Range("H94").Value = 17

$ws.Range("A95").Value = "Jordania"
$ws.Range("B95").Value = 345
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 110
$ws.Range("E95").Value = 230
$ws.Range("F95").Value = 5
$ws.Range("H95").Value = 5

$ws.Range("A96").Value = "Oman"
$ws.Range("B96").Value = 331
$ws.Range("C96").Value = 33
$ws.Range("D96").Value = 61
$ws.Range("E96").Value = 268
$ws.Range("F96").Value = 3
$ws.Range("H96").Value = 2

# --- Row 100: Estado de Palestina ---
$ws.Range("B100").Value = 253
$ws.Range("C100").Value = 16
$ws.Range("E100").Value = 227

# --- Row 119: Isla de Man ---
$ws.Range("B119").Value = 139
$ws.Range("C119").Value = 12
$ws.Range("D119").Value = 55
$ws.Range("E119").Value = 83

# --- Rows 126-128: Trinidad yTobago / Ruanda / Gibraltar reshuffle ---
$ws.Range("A126").Value = "Gibraltar"
$ws.Range("B126").Value = 109
$ws.Range("C126").Value = 6
$ws.Range("D126").Value = 52
$ws.Range("E126").Value = 57
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 0

$ws.Range("A127").Value = "Trinidad yTobago"
$ws.Range("B127").Value = 105
$ws.Range("C127").Value = 1
$ws.Range("D127").Value = 1
$ws.Range("E127").Value = 96
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 8

$ws.Range("A128").Value = "Ruanda"
$ws.Range("B128").Value = 104
$ws.Range("D128").Value = 4
$ws.Range("E128").Value = 100

# --- Row 158: Benin ---
$ws.Range("B158").Value = 23
$ws.Range("C158").Value = 1
$ws.Range("E158").Value = 17

# --- Row 159: Birmania ---
$ws.Range("B159").Value = 22
$ws.Range("C159").Value = 1
$ws.Range("E159").Value = 21
